# Generate Report for Handoff
# Updates status from "In Translation" to "Ready for handoff" and refreshes
# the related handoff timestamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-29-18 16:29:56"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-18 16:29:53"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-18 16:29:56"
